$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-28 Friday" "2024-06-29 Saturday"

Replace-Text "485×6=" "112×4="
Replace-Text "121×9=" "193×4="
Replace-Text "931×6=" "751×7="
Replace-Text "713×3=" "227×7="
Replace-Text "422×9=" "406×7="
Replace-Text "114×4=" "908×9="
Replace-Text "467×9=" "846×9="
Replace-Text "856×7=" "982×7="
Replace-Text "848×9=" "424×2="
Replace-Text "727×5=" "477×5="
Replace-Text "714×9=" "869×4="
Replace-Text "933×5=" "234×2="
Replace-Text "438×7=" "570×7="
Replace-Text "506×6=" "378×7="
Replace-Text "236×2=" "962×7="
Replace-Text "784×2=" "650×2="
Replace-Text "932×6=" "887×2="
Replace-Text "102×2=" "160×8="
Replace-Text "856×2=" "263×5="
Replace-Text "651×6=" "564×2="
Replace-Text "573×8=" "627×3="
Replace-Text "538×2=" "620×9="
Replace-Text "975×8=" "332×6="
Replace-Text "404×7=" "585×4="
Replace-Text "318×2=" "799×5="
